$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 7354878
$ws.Range("I40").Value = 2007.2037
$ws.Range("J40").Value = 35715950
$ws.Range("K40").Value = 2007.2037
$ws.Range("L40").Value = 35715950
$ws.Range("M40").Value = -1832.2037
$ws.Range("N40").Value = -35716300

# Row 53
$ws.Range("H53").Value = 200500.6
$ws.Range("I53").Value = 500150.5
$ws.Range("J53").Value = 734
$ws.Range("K53").Value = 500150.5
$ws.Range("L53").Value = 734
$ws.Range("M53").Value = -499513.5
$ws.Range("N53").Value = -2008

# Row 64
$ws.Range("H64").Value = 5166
$ws.Range("I64").Value = 8580
$ws.Range("J64").Value = 3459
$ws.Range("K64").Value = 8580
$ws.Range("L64").Value = 3459
$ws.Range("M64").Value = -8332
$ws.Range("N64").Value = -3955

# Row 67
$ws.Range("H67").Value = 5166
$ws.Range("I67").Value = 8580
$ws.Range("J67").Value = 3459
$ws.Range("K67").Value = 8580
$ws.Range("L67").Value = 3459
$ws.Range("M67").Value = -7722
$ws.Range("N67").Value = -5175

# Row 74
$ws.Range("H74").Value = 3938.4614
$ws.Range("I74").Value = 4036.3635
$ws.Range("J74").Value = 3400
$ws.Range("K74").Value = 4036.3635
$ws.Range("L74").Value = 3400
$ws.Range("M74").Value = -3100.3635
$ws.Range("N74").Value = -5272

# Row 77
$ws.Range("H77").Value = 3938.4614
$ws.Range("I77").Value = 4036.3635
$ws.Range("J77").Value = 3400
$ws.Range("K77").Value = 20181.8175
$ws.Range("L77").Value = 17000
$ws.Range("M77").Value = -15501.8175
$ws.Range("N77").Value = -26360

# Row 138
$ws.Range("H138").Value = 2442.7576
$ws.Range("I138").Value = 1241.1578
$ws.Range("J138").Value = 3191.2952
$ws.Range("K138").Value = 3723.4734
$ws.Range("L138").Value = 9573.8856
$ws.Range("M138").Value = 1416.5266
$ws.Range("N138").Value = -19853.8856

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 318374.53
$ws.Range("I61").Value = 8013.8237
$ws.Range("J61").Value = 670116.7
$ws.Range("K61").Value = 8013.8237
$ws.Range("L61").Value = 670116.7
$ws.Range("M61").Value = -7801.8237
$ws.Range("N61").Value = -670540.7

# Row 63
$ws.Range("H63").Value = 100011980
$ws.Range("I63").Value = 125014350
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 125014350
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -125013664
$ws.Range("N63").Value = -3872

# Row 66
$ws.Range("H66").Value = 100011980
$ws.Range("I66").Value = 125014350
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 625071750
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -625068318
$ws.Range("N66").Value = -19364

# Row 74
$ws.Range("H74").Value = 1869.4517
$ws.Range("I74").Value = 1583.2
$ws.Range("J74").Value = 2389.9092
$ws.Range("K74").Value = 1583.2
$ws.Range("L74").Value = 2389.9092
$ws.Range("M74").Value = -709.2
$ws.Range("N74").Value = -4137.9092

# Row 77
$ws.Range("H77").Value = 1869.4517
$ws.Range("I77").Value = 1583.2
$ws.Range("J77").Value = 2389.9092
$ws.Range("K77").Value = 7916
$ws.Range("L77").Value = 11949.546
$ws.Range("M77").Value = -3548
$ws.Range("N77").Value = -20685.546

# Row 132
$ws.Range("H132").Value = 1788096.6
$ws.Range("I132").Value = 1492.0889
$ws.Range("J132").Value = 9096933
$ws.Range("K132").Value = 4476.2667
$ws.Range("L132").Value = 27290799
$ws.Range("M132").Value = -1946.2667
$ws.Range("N132").Value = -27295859

# Row 136
$ws.Range("H136").Value = 318374.53
$ws.Range("I136").Value = 8013.8237
$ws.Range("J136").Value = 670116.7
$ws.Range("K136").Value = 24041.4711
$ws.Range("L136").Value = 2010350.1
$ws.Range("M136").Value = -21491.4711
$ws.Range("N136").Value = -2015450.1

# Row 139
$ws.Range("H139").Value = 51699.445
$ws.Range("J139").Value = 51699.445
$ws.Range("L139").Value = 51699.445
$ws.Range("N139").Value = -61979.445

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2183.3333
$ws.Range("I105").Value = 2077.7778
$ws.Range("K105").Value = 2077.7778
$ws.Range("M105").Value = -330.7777999999998

# Row 134
$ws.Range("H134").Value = 22361.203
$ws.Range("I134").Value = 4226.5684
$ws.Range("J134").Value = 102153.6
$ws.Range("K134").Value = 12679.7052
$ws.Range("L134").Value = 306460.8
$ws.Range("M134").Value = -10144.7052
$ws.Range("N134").Value = -311530.8

$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# Row 30
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

# Row 82
$ws.Range("H82").Value = 39890.832
$ws.Range("J82").Value = 41836.2
$ws.Range("L82").Value = 41836.2
$ws.Range("N82").Value = -42558.2

# Row 85
$ws.Range("H85").Value = 39890.832
$ws.Range("J85").Value = 41836.2
$ws.Range("L85").Value = 41836.2
$ws.Range("N85").Value = -44332.2

# Row 88
$ws.Range("H88").Value = 38842.715
$ws.Range("J88").Value = 38842.715
$ws.Range("L88").Value = 38842.715
$ws.Range("N88").Value = -39654.715

# Row 91
$ws.Range("H91").Value = 38842.715
$ws.Range("J91").Value = 38842.715
$ws.Range("L91").Value = 38842.715
$ws.Range("N91").Value = -41650.715

# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# Row 130
$ws.Range("H130").Value = 100000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 100000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 100000
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -110040

# Row 134
$ws.Range("H134").Value = 297628.75
$ws.Range("I134").Value = 3965.0908
$ws.Range("J134").Value = 836012.0600000001
$ws.Range("K134").Value = 11895.2724
$ws.Range("L134").Value = 2508036.18
$ws.Range("M134").Value = -9360.2724
$ws.Range("N134").Value = -2513106.18

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 472.2
$ws.Range("I107").Value = 450.69232
$ws.Range("J107").Value = 512.1429000000001
$ws.Range("K107").Value = 1352.07696
$ws.Range("L107").Value = 1536.4287
$ws.Range("M107").Value = 567.9230400000001
$ws.Range("N107").Value = -5376.4287

# Row 136
$ws.Range("H136").Value = 10099.23
$ws.Range("I136").Value = 10117.272
$ws.Range("K136").Value = 30351.816
$ws.Range("M136").Value = -25251.816

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 11291.363
$ws.Range("I80").Value = 12170.5
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 12170.5
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -11172.5
$ws.Range("N80").Value = -4496

# Row 83
$ws.Range("H83").Value = 11291.363
$ws.Range("I83").Value = 12170.5
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 60852.5
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -55860.5
$ws.Range("N83").Value = -22484

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 10069.594
$ws.Range("I136").Value = 7403.174
$ws.Range("J136").Value = 16883.777
$ws.Range("K136").Value = 22209.522
$ws.Range("L136").Value = 50651.33099999999
$ws.Range("M136").Value = -19659.522
$ws.Range("N136").Value = -55751.33099999999

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 2396.4927
$ws.Range("I136").Value = 2166.45
$ws.Range("K136").Value = 6499.349999999999
$ws.Range("M136").Value = -3949.349999999999
